# Refresh the cryptos price/volume table (and two re-ranked rows) to match
# the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writes $value into $cellRef as literal text. Many values in this sheet look
# like plain numbers (prices such as "535.95" or "5.60") but are stored as
# text in the source data, so a bare .Value assignment would let Excel's COM
# layer coerce them into doubles (silently dropping formatting such as the
# trailing zero in "5.60"). A leading apostrophe reproduces what a user typing
# the value into the grid would do: Excel keeps it as quote-prefixed text.
function Set-TextValue($cellRef, $value) {
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        $ws.Range($cellRef).Value = "'" + $value
    } else {
        $ws.Range($cellRef).Value = $value
    }
}

# Each entry is ,(cellRef, newValue); the leading comma keeps the pair from
# being flattened into the outer @(...) list.
$updates = @(
    ,('D2', '58.967.77')
    ,('E2', '  -4.94%  ')
    ,('D3', '2.469.27')
    ,('E3', '  -4.56%  ')
    ,('E4', '  +0.12%  ')
    ,('D5', '535.95')
    ,('E5', '  -3.23%  ')
    ,('D6', '144.37')
    ,('E6', '  -6.60%  ')
    ,('D7', '0.998')
    ,('E7', '  -0.12%  ')
    ,('D8', '0.569')
    ,('E8', '  -4.05%  ')
    ,('D9', '2.498.49')
    ,('E9', '  -3.80%  ')
    ,('D10', '0.100')
    ,('E10', '  -4.19%  ')
    ,('E11', '  -1.91%  ')
    ,('D12', '5.60')
    ,('E12', '  +2.20%  ')
    ,('D13', '0.355')
    ,('E13', '  -2.85%  ')
    ,('D14', '2.922.97')
    ,('E14', '  -4.09%  ')
    ,('D15', '23.95')
    ,('E15', '  -6.25%  ')
    ,('D16', '58.931.69')
    ,('E16', '  -4.82%  ')
    ,('E17', '  -3.88%  ')
    ,('D18', '2.494.89')
    ,('E18', '  -3.74%  ')
    ,('D19', '11.39')
    ,('E19', '  -2.18%  ')
    ,('D20', '4.33')
    ,('E20', '  -4.85%  ')
    ,('D21', '324.26')
    ,('E21', '  -4.34%  ')
    ,('D22', '0.996')
    ,('E22', '  -0.16%  ')
    ,('D23', '5.76')
    ,('E23', '  -4.61%  ')
    ,('D24', '60.78')
    ,('E24', '  -3.03%  ')
    ,('D25', '0.442')
    ,('E25', '  -11.43%  ')
    ,('B26', 'Binance-PegBSC-USD')
    ,('C26', 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd')
    ,('D26', '0.998')
    ,('E26', '  -0.19%  ')
    ,('D27', '2.609.31')
    ,('E27', '  -3.70%  ')
    ,('B28', 'Kaspa')
    ,('C28', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas')
    ,('D28', '0.161')
    ,('E28', '  -4.04%  ')
    ,('D29', '7.81')
    ,('E29', '  -3.31%  ')
    ,('D30', '7.05')
    ,('E30', '  -0.68%  ')
    ,('D31', '1.28')
    ,('E31', '  -2.92%  ')
    ,('D32', '0.0₃0777')
    ,('E32', '  -7.22%  ')
    ,('D33', '1.80')
    ,('E33', '  -6.36%  ')
    ,('D34', '0.997')
    ,('E34', '  -0.16%  ')
    ,('D35', '158.24')
    ,('E35', '  -1.14%  ')
    ,('D36', '1.41')
    ,('E36', '  -0.78%  ')
    ,('D37', '18.54')
    ,('E37', '  -3.64%  ')
    ,('D38', '4.44')
    ,('E38', '  -5.88%  ')
    ,('D39', '1.64')
    ,('E39', '  -7.65%  ')
    ,('D40', '5.85')
    ,('E40', '  -2.54%  ')
    ,('D41', '310.95')
    ,('E41', '  -8.76%  ')
    ,('D42', '36.68')
    ,('E42', '  -2.21%  ')
    ,('D43', '3.72')
    ,('E43', '  -5.08%  ')
    ,('D44', '0.824')
    ,('E44', '  -7.88%  ')
    ,('D45', '0.997')
    ,('E45', '  -0.09%  ')
    ,('D46', '0.595')
    ,('E46', '  -1.99%  ')
    ,('E47', '  -1.75%  ')
    ,('D48', '124.23')
    ,('E48', '  -0.53%  ')
    ,('B49', 'Stellar')
    ,('C49', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm')
    ,('D49', '0.0929')
    ,('E49', '  -3.71%  ')
    ,('B50', 'Hedera')
    ,('C50', 'https://coinranking.com/coin/jad286TjB+hedera-hbar')
    ,('D50', '0.0523')
    ,('E50', '  -4.52%  ')
    ,('D51', '0.0229')
    ,('E51', '  -4.51%  ')
)

foreach ($u in $updates) {
    Set-TextValue $u[0] $u[1]
}
